$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "FirstSheet" worksheet as the first tab, keep "Sheet1"
#    as the second tab and push the original "CellValue" sheet to third.
# ---------------------------------------------------------------------------
$first = $wb.Worksheets.Add()
$first.Name = "FirstSheet"
$first.Move($wb.Worksheets.Item(1))

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Move($wb.Worksheets.Item(2))

$cellValue = $wb.Worksheets.Item("CellValue")

# ---------------------------------------------------------------------------
# 2. Populate "FirstSheet" with the BDD parameter-table sample data.
# ---------------------------------------------------------------------------

# Reuse already-defined cell styles from the "CellValue" sheet so we do not
# fork duplicate style entries: copy + paste-special(formats) is format-only.
$cellValue.Range("B3").Copy()
$first.Range("A2:E2").PasteSpecial(-4122)
$first.Range("G2:H2").PasteSpecial(-4122)
$first.Range("A3:E3").PasteSpecial(-4122)
$first.Range("G3:H3").PasteSpecial(-4122)

$cellValue.Range("A3").Copy()
$first.Range("C4:C10").PasteSpecial(-4122)

$cellValue.Range("F4").Copy()
$first.Range("D4:E10").PasteSpecial(-4122)
$first.Range("G4:G9").PasteSpecial(-4122)
$first.Range("H4:H9").PasteSpecial(-4122)

$cellValue.Range("K4").Copy()
$first.Range("D10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Gray "locked" look for columns A & B (rows 4-10) - new fill/border combo.
$first.Range("A4:B10").Interior.Color = 14277081
$first.Range("A4:B10").Borders.LineStyle = 1

# "Neutral"(warning) look for the invalid "wrongheader" column (F).
$first.Range("F2:F10").Style = "Neutral"
$first.Range("F3:F10").Borders.LineStyle = 1
$first.Range("F2").Borders.Item(7).LineStyle = 1
$first.Range("F2").Borders.Item(10).LineStyle = 1
$first.Range("F2").Borders.Item(8).LineStyle = -4142
$first.Range("F2").Borders.Item(9).LineStyle = -4142

# Plain gray cell (no border) used for the free-text note below the table.
$first.Range("A13").Interior.Color = 14277081

# H10 numeric value is right aligned (matches the green block elsewhere).
$first.Range("H10").HorizontalAlignment = -4152

# Row 2 - header labels.
$first.Range("A2").Value = "Step"
$first.Range("B2").Value = "Role"
$first.Range("C2").Value = "Parameter Name\ Header"
$first.Range("D2").Value = "Scenario1"
$first.Range("E2").Value = "Scenario2"
$first.Range("F2").Value = "wrongheader"
$first.Range("G2").Value = "Scenario3"
$first.Range("H2").Value = "Scenario4"

# Row 4 - step1.
$first.Range("A4").Value = "step1,blabla"
$first.Range("B4").Value = "user"
$first.Range("C4").Value = "ParamName1"
$first.Range("D4").Value = "V1.1"
$first.Range("E4").Value = "V1.2"
$first.Range("F4").Value = "V1.5"
$first.Range("G4").Value = "V1.3"
$first.Range("H4").Value = "V1.4"

# Row 5.
$first.Range("C5").Value = "NA"

# Row 6.
$first.Range("A6").Value = "bbb"

# Row 8 - step2.
$first.Range("A8").Value = "step2,blabla"
$first.Range("B8").Value = "system"
$first.Range("C8").Value = "ParamName2"
$first.Range("D8").Value = "V2.1"
$first.Range("F8").Value = "V2.5"
$first.Range("G8").Value = "V2.3"
$first.Range("H8").Value = "V2.4"

# Row 9 - step3.
$first.Range("A9").Value = "step3,blabla"
$first.Range("B9").Value = "user"
$first.Range("C9").Value = "ParamName3"

# Row 10 - step4.
$first.Range("A10").Value = "step4,blabla"
$first.Range("B10").Value = "system"
$first.Range("C10").Value = "ParamName4"
$first.Range("D10").NumberFormat = "@"
$first.Range("D10").Value = "2021/4/30"
$first.Range("E10").Value = $false
$first.Range("F10").Value = "V4.5"
$first.Range("G10").Value = $true
$first.Range("H10").Value = 4.4

# Row 13 - free note.
$first.Range("A13").Value = "eee"

# Rows 15/16 - helper cells referenced by the E8 formula.
$first.Range("C15").Value = "useless"
$first.Range("D15").Value = "no"
$first.Range("E15").Value = "V2.2"
$first.Range("G15").Value = "no"
$first.Range("H15").Value = "no"
$first.Range("F16").Value = "no"

# E8 looks the current value up from E15 via formula.
$first.Range("E8").Formula = "=E15"

$first.Columns.Item(1).AutoFit()
$first.Columns.Item(2).AutoFit()
$first.Columns.Item(3).AutoFit()
$first.Columns.Item(4).AutoFit()
$first.Columns.Item(6).AutoFit()

$first.Range("A1").Select()
